{"js": "// 1. Update the \"Logical Expression\" summary line (unique full-line text, safest to search/replace whole).\nconst body = context.document.body;\nconst titleResults = body.search(\"Logical Expression: a,a\u2192v\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\n    \"Logical Expression: p\u2192(q\u2192r),p\u2192q,p\u22a2r\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2. Locate the results table and update the two existing \"Expression\" cells.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Row 1 (0-based, i.e. the first data row \"Line 1\"), column 1 (\"Expression\"): \"a\" -> \"p\"\nconst cellRow1 = table.getCell(1, 1);\nconst cellRow1Results = cellRow1.body.search(\"a\", { matchCase: true });\ncellRow1Results.load(\"items\");\nawait context.sync();\ncellRow1Results.items[0].insertText(\"p\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Row 2 (0-based, i.e. the second data row \"Line 2\"), column 1 (\"Expression\"): \"a\u2192v\" -> \"p\u2192(q\u2192r)\"\nconst cellRow2 = table.getCell(2, 1);\nconst cellRow2Results = cellRow2.body.search(\"a\u2192v\", { matchCase: true });\ncellRow2Results.load(\"items\");\nawait context.sync();\ncellRow2Results.items[0].insertText(\"p\u2192(q\u2192r)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Append three new proof lines continuing the Modus Ponens derivation.\ntable.addRows(Word.InsertLocation.end, 3, [\n  [\"3\", \"q\u2192r\", \"MP\", \"1\", \"2\", \"\"],\n  [\"4\", \"p\u2192q\", \"Data\", \"\", \"\", \"\"],\n  [\"5\", \"q\", \"MP\", \"1\", \"4\", \"\"]\n]);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the \"Logical Expression\" summary line.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"Logical Expression: a,a\u2192v\", $false, $false, $false, $false, $false, $true, 1, $false, \"Logical Expression: p\u2192(q\u2192r),p\u2192q,p\u22a2r\", 2) | Out-Null\n\n# 2. Update the two existing \"Expression\" cells in the results table.\n$t = $d.Tables.Item(1)\n$t.Cell(2,2).Range.Text = \"p\"\n$t.Cell(3,2).Range.Text = \"p\u2192(q\u2192r)\"\n\n# 3. Append three new proof lines continuing the Modus Ponens derivation:\n#    3: q\u2192r   MP   1 2\n#    4: p\u2192q   Data\n#    5: q     MP   1 4\n$row3 = $t.Rows.Add()\n$row3.Cells(1).Range.Text = \"3\"\n$row3.Cells(2).Range.Text = \"q\u2192r\"\n$row3.Cells(3).Range.Text = \"MP\"\n$row3.Cells(4).Range.Text = \"1\"\n$row3.Cells(5).Range.Text = \"2\"\n\n$row4 = $t.Rows.Add()\n$row4.Cells(1).Range.Text = \"4\"\n$row4.Cells(2).Range.Text = \"p\u2192q\"\n$row4.Cells(3).Range.Text = \"Data\"\n\n$row5 = $t.Rows.Add()\n$row5.Cells(1).Range.Text = \"5\"\n$row5.Cells(2).Range.Text = \"q\"\n$row5.Cells(3).Range.Text = \"MP\"\n$row5.Cells(4).Range.Text = \"1\"\n$row5.Cells(5).Range.Text = \"4\"\n"}
